$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 495
$ws1.Range("F6").Value = 679

# Sheet "全部类型"
$ws2 = $wb.Worksheets.Item("全部类型")
$ws2.Range("F3").Value = 495
$ws2.Range("F6").Value = 679
